$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.022.07"
$ws.Range("E2").Value = "  -0.92%  "
$ws.Range("D3").Value = "3.579.30"
$ws.Range("E3").Value = "  -1.85%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'580.19"
$ws.Range("E5").Value = "  -2.82%  "
$ws.Range("D6").Value = "'188.54"
$ws.Range("E6").Value = "  -2.66%  "
$ws.Range("E7").Value = "  -2.90%  "
$ws.Range("D8").Value = "3.577.19"
$ws.Range("E8").Value = "  -0.78%  "
$ws.Range("D9").Value = "'1.00"
$ws.Range("E9").Value = "  -0.02%  "
$ws.Range("D10").Value = "'0.178"
$ws.Range("E10").Value = "  -3.08%  "
$ws.Range("D11").Value = "'0.659"
$ws.Range("E11").Value = "  -1.20%  "
$ws.Range("D12").Value = "'56.02"
$ws.Range("E12").Value = "  -3.69%  "
$ws.Range("D13").Value = "'0.0000300"
$ws.Range("E13").Value = "  -1.17%  "
$ws.Range("D14").Value = "'9.68"
$ws.Range("E14").Value = "  -1.36%  "
$ws.Range("D15").Value = "4.148.06"
$ws.Range("E15").Value = "  -1.67%  "
$ws.Range("D16").Value = "'19.96"
$ws.Range("E16").Value = "  +2.51%  "
$ws.Range("D17").Value = "3.576.13"
$ws.Range("E17").Value = "  -1.66%  "
$ws.Range("D18").Value = "69.942.39"
$ws.Range("E18").Value = "  -0.88%  "
$ws.Range("D19").Value = "'12.56"
$ws.Range("E19").Value = "  -1.09%  "
$ws.Range("E20").Value = "  +0.12%  "
$ws.Range("E21").Value = "  -1.95%  "
$ws.Range("D22").Value = "'473.79"
$ws.Range("E22").Value = "  -5.24%  "
$ws.Range("D23").Value = "'19.21"
$ws.Range("E23").Value = "  +14.21%  "
$ws.Range("E24").Value = "  -8.47%  "
$ws.Range("D25").Value = "'4.35"
$ws.Range("E25").Value = "  -2.76%  "
$ws.Range("D26").Value = "'88.54"
$ws.Range("E26").Value = "  -3.18%  "
$ws.Range("E27").Value = "  -2.53%  "
$ws.Range("D28").Value = "'11.01"
$ws.Range("E28").Value = "  -2.92%  "
$ws.Range("D29").Value = "'9.36"
$ws.Range("E29").Value = "  -0.85%  "
$ws.Range("D30").Value = "'32.13"
$ws.Range("E30").Value = "  -1.62%  "
$ws.Range("D31").Value = "'7.64"
$ws.Range("E31").Value = "  +1.02%  "
$ws.Range("D32").Value = "'0.120"
$ws.Range("E32").Value = "  +2.46%  "
$ws.Range("D33").Value = "'12.09"
$ws.Range("E33").Value = "  -1.25%  "
$ws.Range("D34").Value = "'65.77"
$ws.Range("E34").Value = "  -0.08%  "
$ws.Range("D35").Value = "'581.32"
$ws.Range("E35").Value = "  -6.12%  "
$ws.Range("D36").Value = "'38.76"
$ws.Range("E36").Value = "  +1.30%  "
$ws.Range("E37").Value = "  -0.02%  "
$ws.Range("D38").Value = "0.0₃0800"
$ws.Range("E38").Value = "  -4.81%  "
$ws.Range("D39").Value = "'0.395"
$ws.Range("E39").Value = "  -2.14%  "
$ws.Range("E40").Value = "  -6.06%  "
$ws.Range("E41").Value = "  +14.69%  "
$ws.Range("D42").Value = "'3.51"
$ws.Range("E42").Value = "  -5.96%  "
$ws.Range("B43").Value = "Fetch.AI"
$ws.Range("C43").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D43").Value = "'2.87"
$ws.Range("E43").Value = "  +6.59%  "
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "3.227.64"
$ws.Range("E44").Value = "  -3.44%  "
$ws.Range("D45").Value = "'3.12"
$ws.Range("E45").Value = "  +0.42%  "
$ws.Range("E46").Value = "  -1.64%  "
$ws.Range("D47").Value = "'9.54"
$ws.Range("E47").Value = "  +4.61%  "
$ws.Range("E48").Value = "  +1.01%  "
$ws.Range("E49").Value = "  -0.61%  "
$ws.Range("D50").Value = "'0.998"
$ws.Range("E50").Value = "  -0.06%  "
$ws.Range("E51").Value = "  -4.10%  "
